# Update "Alert Parameters (working).xlsx" - Third Iteration sheet
# - "Power" -> "POWER" and "Temp" -> "TEMP" in the metric column (B)
# - Selection moved from F21 to B42 (and scroll position nudged accordingly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose metric (column B) was the literal "Power" -> now "POWER"
$powerRows = 4,5,6,7,8,13,14,16,19,22,25
foreach ($r in $powerRows) {
    $ws.Cells.Item($r, 2).Value = "POWER"
}

# Rows whose metric (column B) was the literal "Temp" -> now "TEMP"
$tempRows = 9,10,11,12,15,18,21,24,27,28,29,30,31,32,33,34,35,36,37,38
foreach ($r in $tempRows) {
    $ws.Cells.Item($r, 2).Value = "TEMP"
}

# Scroll the view down a bit and move the active selection to B42
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B42").Select()
